# Add Olin CoLab Location
# Inserts a new row above the existing row 90 ("Library Tech Services")
# for a new "Olin Library Digital CoLab" location mapping, pushing all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a blank row at row 90, shifting row 90 (and everything below it)
# down to row 91.
$ws.Rows.Item(90).Insert()

# Populate the new row: Voyager Display Name (A) and Facet Display Name (E).
# Columns B (Suppress), C (Call#), D (Holding Note) are left blank, matching
# the other simple single-mapping rows in the sheet.
$ws.Range("A90").Value = "Olin Library Digital CoLab"
$ws.Range("E90").Value = "Olin Library > Digital CoLab Room 701"

# Reflect the author's cursor ending up on the newly-added row.
$ws.Range("A90").Select()
